$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows before row 333 (existing rows 333-382 shift down to 340-389).
$ws.Rows("333:339").Insert()

# Static column values shared by every data row in this block (Mercado ID ..
# Categoría, and Clasificación at the end never change within this sheet).
$mercadoId = 6
$mercado   = 'Mercado Mayorista Lo Valledor de Santiago'
$region    = 'Metropolitana'
$codreg    = 13
$catId     = 100112013
$categoria = 'Alcachofa'
$clasif    = 'Hortaliza'

function Set-Row($r, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidadComerc, $origen, $precioKg, $kgUnidades) {
    $ws.Cells.Item($r, 1).Value2  = $mercadoId
    $ws.Cells.Item($r, 2).Value2  = $mercado
    $ws.Cells.Item($r, 3).Value2  = $region
    $ws.Cells.Item($r, 4).Value2  = $fecha
    $ws.Cells.Item($r, 5).Value2  = $codreg
    $ws.Cells.Item($r, 6).Value2  = $catId
    $ws.Cells.Item($r, 7).Value2  = $categoria
    $ws.Cells.Item($r, 8).Value2  = $variedad
    $ws.Cells.Item($r, 9).Value2  = $calidad
    $ws.Cells.Item($r, 10).Value2 = $volumen
    $ws.Cells.Item($r, 11).Value2 = $precioMin
    $ws.Cells.Item($r, 12).Value2 = $precioMax
    $ws.Cells.Item($r, 13).Value2 = $precioProm
    $ws.Cells.Item($r, 14).Value2 = $unidadComerc
    $ws.Cells.Item($r, 15).Value2 = $origen
    $ws.Cells.Item($r, 16).Value2 = $precioKg
    $ws.Cells.Item($r, 17).Value2 = $kgUnidades
    $ws.Cells.Item($r, 18).Value2 = $clasif
}

Set-Row 333 44474 'Argentina(o)' 'Primera' 590  12000 14000 12915 '$/caja 50 unidades' 'Provincia de Limarí'   258   50
Set-Row 334 44474 'Española'     'Extra'   580  12000 13000 12448 '$/caja 25 unidades' 'Provincia de Limarí'   12448 1
Set-Row 335 44474 'Española'     'Primera' 430  10000 12000 11209 '$/caja 30 unidades' 'Provincia de Limarí'   374   30
Set-Row 336 44474 'Española'     'Primera' 7500 500   550   517   '$/unidad'           'Región Metropolitana'  517   1
Set-Row 337 44474 'Española'     'Segunda' 5000 400   450   420   '$/unidad'           'Región Metropolitana'  420   1
Set-Row 338 44474 'Española'     'Tercera' 3500 300   350   321   '$/unidad'           'Región Metropolitana'  321   1
Set-Row 339 44474 'Madrigal'     'Primera' 680  8000  10000 8941  '$/caja 40 unidades' 'Provincia de Limarí'   224   40

# Keep the date cells in column D using the same datetime display format as
# the rest of that column.
$dateFormat = $ws.Range("D340").NumberFormat
$ws.Range("D333:D339").NumberFormat = $dateFormat
